$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Date value (row 8, column B)
$ws.Range("B8").Value = "2024-06-26T18:33:17-04:00"

# Correct the Context values to remove the StructureDefinition URL prefixes
$ws.Range("B21").Value = "element:Patient"
$ws.Range("B22").Value = "element:Practitioner"
